$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F14").Value = 891
$ws1.Range("F15").Value = 879
$ws1.Range("F19").Value = 663
$ws1.Range("F20").Value = 780
$ws1.Range("F22").Value = 2936
$ws1.Range("F29").Value = 254
$ws1.Range("F37").Value = 1779
$ws1.Range("F41").Value = 194

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 891
$ws4.Range("F13").Value = 879
$ws4.Range("F19").Value = 780
$ws4.Range("F21").Value = 2936
$ws4.Range("F27").Value = 601
$ws4.Range("F28").Value = 254
$ws4.Range("F43").Value = 1779
$ws4.Range("F46").Value = 194
